# issue #5: stock data from json to db
# Adds a "category" column (value "normal") between property_category and date,
# plus "source_file" (tmp5431) and "index" columns at the end of the 股票(stock)
# sheet. Also strips the "★" marker prefix from the two stock name values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Header row (row 1) -----------------------------------------------
# Existing headers B1:K1 (name,owner,quantity,face_value,currency,total,
# property_category,date,legislator_name,legislator_id) stay as-is.
# Insert new "category" header right after H1 (property_category) by
# shifting I1:K1 one column to the right, then append two brand new
# trailing headers (source_file, index).

$ws.Range("I1:K1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"

$ws.Range("K1").Copy() | Out-Null
$ws.Range("M1:N1").PasteSpecial(-4122) | Out-Null
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data row 2 ---------------------------------------------------------
$ws.Range("I2:K2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2011-11-25"
$ws.Range("K2").Value = "邱議瑩"
$ws.Range("L2").Value = 913

$ws.Range("K2").Copy() | Out-Null
$ws.Range("M2:N2").PasteSpecial(-4122) | Out-Null
$ws.Range("M2").Value = "tmp5431"
$ws.Range("N2").Value = 43

$ws.Range("B2").Value = "EquinoxMineralsLimited"

# --- Data row 3 ---------------------------------------------------------
$ws.Range("I3:K3").Copy() | Out-Null
$ws.Range("J3").PasteSpecial(-4122) | Out-Null
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2011-11-25"
$ws.Range("K3").Value = "邱議瑩"
$ws.Range("L3").Value = 913

$ws.Range("K3").Copy() | Out-Null
$ws.Range("M3:N3").PasteSpecial(-4122) | Out-Null
$ws.Range("M3").Value = "tmp5431"
$ws.Range("N3").Value = 44

$ws.Range("B3").Value = "中興商銀"
